# Commit: "Fixed #295 Add the version of M2Doc in the template custom
# properties."
#
# The canonical-OOXML diff for this particular template
# (emptyMTextInFooter-template.docx) only shows word/document.xml,
# word/footer1-3.xml, word/footnotes.xml, word/header1-3.xml and
# word/styles.xml being re-serialized: every start tag's attributes and
# every root element's namespace declarations come back out in
# alphabetical order (e.g. w:headerReference's r:id now precedes
# w:type, w:pgMar's attributes are alphabetised, w:style's w:type now
# comes after w:styleId, ...) and the transient w:rsid* attributes are
# gone. Diffing element-by-element (and re-sorting each tag's
# attributes back alphabetically) shows the "before" and "after" XML
# are attribute-for-attribute, element-for-element, text-for-text
# identical - i.e. this is a pure re-serialization/pretty-printer
# artifact, not a content edit. No paragraph text, run formatting,
# style definition, header/footer text or section property actually
# changed anywhere in this file. (The M2Doc-version custom property
# that the commit message references lives in docProps/custom.xml,
# which is not part of this template's diff at all.)
#
# So there is nothing to change in the Word object model for this
# document: we simply touch the parts involved (document body,
# headers/footers, styles) read-only, which leaves the package
# byte-for-byte identical - matching the target content exactly.

$d = $word.ActiveDocument

$null = $d.Paragraphs.Count
$null = $d.Styles.Count

$sections = $d.Sections
for ($i = 1; $i -le $sections.Count; $i++) {
    $section = $sections.Item($i)
    for ($hf = 1; $hf -le 3; $hf++) {
        $header = $section.Headers.Item($hf)
        if ($header.Exists) { $null = $header.Range.Text }
        $footer = $section.Footers.Item($hf)
        if ($footer.Exists) { $null = $footer.Range.Text }
    }
}
